# Add Selenium web driver support:
#  - Rename element types (typo/cleanup): TestBox -> TextBox, btn -> Button, div -> Div
#  - Move the "ID" value on the Elements sheet from column D (Name) to column C (ID)
#  - Add a new "Assert" execution step (row 4 / ID 4) on the ExecData sheet
#  - Fix up the WaitPeriod for the "Click" step
#  - Update selections / active sheet to the Elements sheet

$wb = $excel.ActiveWorkbook

$execData = $wb.Worksheets.Item("ExecData")
$elements = $wb.Worksheets.Item("Elements")

# --- Elements sheet: fix up element type labels and move the ID value ---
$elements.Range("C2").Value = "lst-ib"
$elements.Range("D2").ClearContents()

$elements.Range("H2").Value = "TextBox"
$elements.Range("H3").Value = "Button"
$elements.Range("H4").Value = "Div"

# --- ExecData sheet: fix WaitPeriod for step 3 (Click), add new step 4 (Assert) ---
$execData.Range("G4").Value = 0

$execData.Range("A5").Value = 4
$execData.Range("B5").Value = "Search Google"
$execData.Range("C5").Value = "Assert"
$execData.Range("G5").Value = 2

# --- Selections: ExecData -> C3, Elements -> C2 (Elements ends up the active tab) ---
$execData.Range("C3").Select()
$elements.Range("C2").Select()
